$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (avoids numeric
# auto-conversion / floating point drift for numeric-looking strings),
# then restore the cell to its original (default) style so no stray
# number-format attribute is left behind on the cell.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "42.289.51"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "2.275.17"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "308.89"
$ws.Range("E5").Value = "  +0.34%  "
Set-TextValue $ws.Range("D6") "97.32"
$ws.Range("E6").Value = "  -0.17%  "
Set-TextValue $ws.Range("D7") "0.528"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.69%  "
Set-TextValue $ws.Range("D10") "35.21"
$ws.Range("E10").Value = "  -3.17%  "
Set-TextValue $ws.Range("D11") "0.0805"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("E12").Value = "  +0.47%  "
Set-TextValue $ws.Range("D13") "6.78"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "2.627.77"
$ws.Range("E14").Value = "  -0.85%  "
Set-TextValue $ws.Range("D15") "14.70"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "2.275.38"
$ws.Range("E16").Value = "  -1.06%  "
Set-TextValue $ws.Range("D17") "0.792"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "42.132.92"
$ws.Range("E18").Value = "  -0.88%  "
Set-TextValue $ws.Range("D19") "12.30"
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -1.86%  "
Set-TextValue $ws.Range("D21") "5.98"
$ws.Range("E21").Value = "  -0.96%  "
Set-TextValue $ws.Range("D22") "67.74"
$ws.Range("E22").Value = "  -0.45%  "
Set-TextValue $ws.Range("D23") "236.48"
$ws.Range("E23").Value = "  -2.80%  "
Set-TextValue $ws.Range("D24") "2.60"
$ws.Range("E24").Value = "  -0.65%  "
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  -0.06%  "
Set-TextValue $ws.Range("D27") "23.59"
$ws.Range("E27").Value = "  -1.99%  "
Set-TextValue $ws.Range("D28") "37.11"
$ws.Range("E28").Value = "  -0.36%  "
Set-TextValue $ws.Range("D29") "9.54"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("E30").Value = "  +0.25%  "
Set-TextValue $ws.Range("D31") "164.24"
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("E34").Value = "  -0.63%  "
Set-TextValue $ws.Range("D35") "0.0737"
$ws.Range("E35").Value = "  -2.24%  "
Set-TextValue $ws.Range("D36") "17.54"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("E38").Value = "  -3.75%  "
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("E40").Value = "  -0.90%  "
Set-TextValue $ws.Range("D41") "4.13"
$ws.Range("E41").Value = "  -1.40%  "
$ws.Range("E42").Value = "  -6.46%  "
$ws.Range("D43").Value = "1.952.71"
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D44") "0.0282"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D45") "18.84"
$ws.Range("E45").Value = "  -2.52%  "
$ws.Range("E46").Value = "  -4.25%  "
Set-TextValue $ws.Range("D47") "9.79"
$ws.Range("E47").Value = "  -4.84%  "
Set-TextValue $ws.Range("D48") "53.68"
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "2.497.50"
$ws.Range("E49").Value = "  -0.84%  "
Set-TextValue $ws.Range("D50") "92.34"
$ws.Range("E50").Value = "  -0.21%  "
Set-TextValue $ws.Range("D51") "71.54"
$ws.Range("E51").Value = "  -1.87%  "

Write-Host "Applied all crypto list updates"
